# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# sheets to the freshly scraped values (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 10101   # was 10097
$ws1.Range("F20").Value = 7989    # was 7985
$ws1.Range("F21").Value = 12592   # was 12589
$ws1.Range("F24").Value = 15      # was 14
$ws1.Range("F29").Value = 411     # was 300
$ws1.Range("F30").Value = 2814    # was 2813
$ws1.Range("F31").Value = 258     # was 257
$ws1.Range("F33").Value = 7930    # was 7927
$ws1.Range("F34").Value = 1467    # was 1466

# --- Sheet "全部类型" (all categories, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 10101   # was 10097
$ws4.Range("F22").Value = 7989    # was 7986
$ws4.Range("F23").Value = 12592   # was 12589
$ws4.Range("F26").Value = 15      # was 14
$ws4.Range("F32").Value = 2814    # was 2813
$ws4.Range("F34").Value = 258     # was 257
$ws4.Range("F36").Value = 7930    # was 7927

Write-Output "Updated F-column interest counts on 展览 and 全部类型 sheets."
